$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values/formatting of the three requirement rows (23, 24, 25)
$row23B = $ws.Range("B23").Value2
$row24B = $ws.Range("B24").Value2
$row25B = $ws.Range("B25").Value2

$row23C = $ws.Range("C23").Value2
$row24C = $ws.Range("C24").Value2
$row25C = $ws.Range("C25").Value2

# Move LOM3229 entry (currently row 23) to the end (row 25),
# shifting LOB1021 (row 24) up to row 23 and LOM3016 (row 25) up to row 24.
$ws.Range("B23").Value2 = $row24B
$ws.Range("C23").Value2 = $row24C

$ws.Range("B24").Value2 = $row25B
$ws.Range("C24").Value2 = $row25C

$ws.Range("B25").Value2 = $row23B
$ws.Range("C25").Value2 = $row23C
